$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the old "Integral" column (C) to hold the
# variance (STD^2) values. This shifts the old Integral column C -> D
# and the old Time column D -> E, along with all their formulas.
$ws.Range("C1").EntireColumn.Insert()

# Header for the new column
$ws.Range("C1").Value = "Varaince"

# Variance = STD^2 for each sample row
$ws.Range("C2:C11").Formula = "=B2^2"

# Make the new column readable
$null = $ws.Range("C1:C11").EntireColumn.AutoFit()

# Row 13 ("Avg"): recompute average of the new variance column
$ws.Range("C13").Formula = "=AVERAGE(C2:C11)"

# Row 14 ("STD"): we no longer want the STD-of-STD value in column B,
# only the STD of the Integral/Time columns (now D and E) remain.
$ws.Range("B14").ClearContents()

# New row 15: RMS computed as the square root of the averaged variance
$ws.Range("A15").Value = "RMS"
$ws.Range("B15").Formula = "=SQRT(C13)"

# Match the final selection left behind in the saved file
$null = $ws.Range("B16").Select()
